$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 6
$ws.Range("B4").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C4").Value = "Metropolitana"
$ws.Range("D4").Value = 44467
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = 300000000
$ws.Range("G4").Value = "Espárragos"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Banquete"
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = 2000
$ws.Range("N4").Value = "`$/caja 10 kilos"
$ws.Range("O4").Value = "Provincia de Linares"
$ws.Range("P4").Value = 200
$ws.Range("Q4").Value = 10
$ws.Range("R4").Value = "Hortaliza"

# Row 5
$ws.Range("A5").Value = 6
$ws.Range("B5").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C5").Value = "Metropolitana"
$ws.Range("D5").Value = 44467
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = 300000000
$ws.Range("G5").Value = "Espárragos"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 1500
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = 1500
$ws.Range("N5").Value = "`$/caja 10 kilos"
$ws.Range("O5").Value = "Provincia de Linares"
$ws.Range("P5").Value = 150
$ws.Range("Q5").Value = 10
$ws.Range("R5").Value = "Hortaliza"

# Row 6
$ws.Range("A6").Value = 6
$ws.Range("B6").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C6").Value = "Metropolitana"
$ws.Range("D6").Value = 44467
$ws.Range("E6").Value = 13
$ws.Range("F6").Value = 300000000
$ws.Range("G6").Value = "Espárragos"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 1300
$ws.Range("L6").Value = 1300
$ws.Range("M6").Value = 1300
$ws.Range("N6").Value = "`$/caja 10 kilos"
$ws.Range("O6").Value = "Provincia de Linares"
$ws.Range("P6").Value = 130
$ws.Range("Q6").Value = 10
$ws.Range("R6").Value = "Hortaliza"

# Row 7
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C7").Value = "Metropolitana"
$ws.Range("D7").Value = 44161
$ws.Range("E7").Value = 13
$ws.Range("F7").Value = 300000000
$ws.Range("G7").Value = "Espárragos"
$ws.Range("H7").Value = "Verde"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 4300
$ws.Range("K7").Value = 1000
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = 1000
$ws.Range("N7").Value = "`$/kilo"
$ws.Range("O7").Value = "Provincia de Linares"
$ws.Range("P7").Value = 1000
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = "Hortaliza"

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C8").Value = "Metropolitana"
$ws.Range("D8").Value = 44161
$ws.Range("E8").Value = 13
$ws.Range("F8").Value = 300000000
$ws.Range("G8").Value = "Espárragos"
$ws.Range("H8").Value = "Verde"
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 2500
$ws.Range("K8").Value = 800
$ws.Range("L8").Value = 800
$ws.Range("M8").Value = 800
$ws.Range("N8").Value = "`$/kilo"
$ws.Range("O8").Value = "Provincia de Linares"
$ws.Range("P8").Value = 800
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = "Hortaliza"

# Row 9
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C9").Value = "Metropolitana"
$ws.Range("D9").Value = 44167
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = 300000000
$ws.Range("G9").Value = "Espárragos"
$ws.Range("H9").Value = "Verde"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 250
$ws.Range("K9").Value = 13000
$ws.Range("L9").Value = 13000
$ws.Range("M9").Value = 13000
$ws.Range("N9").Value = "`$/bandeja 10 kilos"
$ws.Range("O9").Value = "Provincia de Linares"
$ws.Range("P9").Value = 1300
$ws.Range("Q9").Value = 10
$ws.Range("R9").Value = "Hortaliza"

# Row 10
$ws.Range("A10").Value = 6
$ws.Range("B10").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C10").Value = "Metropolitana"
$ws.Range("D10").Value = 44167
$ws.Range("E10").Value = 13
$ws.Range("F10").Value = 300000000
$ws.Range("G10").Value = "Espárragos"
$ws.Range("H10").Value = "Verde"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 160
$ws.Range("K10").Value = 11000
$ws.Range("L10").Value = 11000
$ws.Range("M10").Value = 11000
$ws.Range("N10").Value = "`$/caja 10 kilos"
$ws.Range("O10").Value = "Provincia de Linares"
$ws.Range("P10").Value = 1100
$ws.Range("Q10").Value = 10
$ws.Range("R10").Value = "Hortaliza"

# Row 11
$ws.Range("A11").Value = 6
$ws.Range("B11").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C11").Value = "Metropolitana"
$ws.Range("D11").Value = 44167
$ws.Range("E11").Value = 13
$ws.Range("F11").Value = 300000000
$ws.Range("G11").Value = "Espárragos"
$ws.Range("H11").Value = "Verde"
$ws.Range("I11").Value = "Segunda"
$ws.Range("J11").Value = 120
$ws.Range("K11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 10000
$ws.Range("N11").Value = "`$/bandeja 10 kilos"
$ws.Range("O11").Value = "Provincia de Linares"
$ws.Range("P11").Value = 1000
$ws.Range("Q11").Value = 10
$ws.Range("R11").Value = "Hortaliza"

# Row 12
$ws.Range("A12").Value = 6
$ws.Range("B12").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C12").Value = "Metropolitana"
$ws.Range("D12").Value = 44160
$ws.Range("E12").Value = 13
$ws.Range("F12").Value = 300000000
$ws.Range("G12").Value = "Espárragos"
$ws.Range("H12").Value = "Verde"
$ws.Range("I12").Value = "Banquete"
$ws.Range("J12").Value = 210
$ws.Range("K12").Value = 13000
$ws.Range("L12").Value = 13000
$ws.Range("M12").Value = 13000
$ws.Range("N12").Value = "`$/bandeja 10 kilos"
$ws.Range("O12").Value = "Región Metropolitana"
$ws.Range("P12").Value = 1300
$ws.Range("Q12").Value = 10
$ws.Range("R12").Value = "Hortaliza"

# Row 13
$ws.Range("A13").Value = 6
$ws.Range("B13").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C13").Value = "Metropolitana"
$ws.Range("D13").Value = 44160
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = 300000000
$ws.Range("G13").Value = "Espárragos"
$ws.Range("H13").Value = "Verde"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 340
$ws.Range("K13").Value = 11000
$ws.Range("L13").Value = 11000
$ws.Range("M13").Value = 11000
$ws.Range("N13").Value = "`$/bandeja 10 kilos"
$ws.Range("O13").Value = "Región Metropolitana"
$ws.Range("P13").Value = 1100
$ws.Range("Q13").Value = 10
$ws.Range("R13").Value = "Hortaliza"

# Row 14
$ws.Range("A14").Value = 6
$ws.Range("B14").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C14").Value = "Metropolitana"
$ws.Range("D14").Value = 44160
$ws.Range("E14").Value = 13
$ws.Range("F14").Value = 300000000
$ws.Range("G14").Value = "Espárragos"
$ws.Range("H14").Value = "Verde"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 4300
$ws.Range("K14").Value = 1200
$ws.Range("L14").Value = 1200
$ws.Range("M14").Value = 1200
$ws.Range("N14").Value = "`$/kilo"
$ws.Range("O14").Value = "Región Metropolitana"
$ws.Range("P14").Value = 1200
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = "Hortaliza"

# Row 15
$ws.Range("A15").Value = 6
$ws.Range("B15").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C15").Value = "Metropolitana"
$ws.Range("D15").Value = 44160
$ws.Range("E15").Value = 13
$ws.Range("F15").Value = 300000000
$ws.Range("G15").Value = "Espárragos"
$ws.Range("H15").Value = "Verde"
$ws.Range("I15").Value = "Segunda"
$ws.Range("J15").Value = 250
$ws.Range("K15").Value = 9000
$ws.Range("L15").Value = 9000
$ws.Range("M15").Value = 9000
$ws.Range("N15").Value = "`$/bandeja 10 kilos"
$ws.Range("O15").Value = "Región Metropolitana"
$ws.Range("P15").Value = 900
$ws.Range("Q15").Value = 10
$ws.Range("R15").Value = "Hortaliza"

# Row 16
$ws.Range("A16").Value = 6
$ws.Range("B16").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C16").Value = "Metropolitana"
$ws.Range("D16").Value = 44160
$ws.Range("D16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E16").Value = 13
$ws.Range("F16").Value = 300000000
$ws.Range("G16").Value = "Espárragos"
$ws.Range("H16").Value = "Verde"
$ws.Range("I16").Value = "Segunda"
$ws.Range("J16").Value = 2500
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = 1000
$ws.Range("N16").Value = "`$/kilo"
$ws.Range("O16").Value = "Región Metropolitana"
$ws.Range("P16").Value = 1000
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = "Hortaliza"

# Row 17
$ws.Range("A17").Value = 6
$ws.Range("B17").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C17").Value = "Metropolitana"
$ws.Range("D17").Value = 44162
$ws.Range("D17").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E17").Value = 13
$ws.Range("F17").Value = 300000000
$ws.Range("G17").Value = "Espárragos"
$ws.Range("H17").Value = "Verde"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 5200
$ws.Range("K17").Value = 1100
$ws.Range("L17").Value = 1100
$ws.Range("M17").Value = 1100
$ws.Range("N17").Value = "`$/kilo"
$ws.Range("O17").Value = "Provincia de Linares"
$ws.Range("P17").Value = 1100
$ws.Range("Q17").Value = 1
$ws.Range("R17").Value = "Hortaliza"

# Row 18
$ws.Range("A18").Value = 6
$ws.Range("B18").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C18").Value = "Metropolitana"
$ws.Range("D18").Value = 44162
$ws.Range("D18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E18").Value = 13
$ws.Range("F18").Value = 300000000
$ws.Range("G18").Value = "Espárragos"
$ws.Range("H18").Value = "Verde"
$ws.Range("I18").Value = "Segunda"
$ws.Range("J18").Value = 3400
$ws.Range("K18").Value = 900
$ws.Range("L18").Value = 900
$ws.Range("M18").Value = 900
$ws.Range("N18").Value = "`$/kilo"
$ws.Range("O18").Value = "Provincia de Linares"
$ws.Range("P18").Value = 900
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = "Hortaliza"
